$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row "back of a napkin" relay setting updates.
# Columns: G (50P IOC Trip Pickup), K (51P TOC Trip Pickup), N (27P Trip Pickup),
#          O (59P Trip Pickup), R (51P TOC Trip Pickup, right block)
$rows = @(
    @{ Row = 2;  G = 1200; K = 0.60140653040586012; N = 0.8; O = 1.2; R = 0.60140653040586012 },
    @{ Row = 3;  G = 1200; K = 2.9124092514083011;  N = 0.8; O = 1.2; R = 2.9124092514083011 },
    @{ Row = 4;  G = 1200; K = 1.3030474825460301;  N = 0.8; O = 1.2; R = 1.3030474825460301 },
    @{ Row = 5;  G = 1200; K = 0.10023442173431002; N = 0.8; O = 1.2; R = 0.10023442173431002 },
    @{ Row = 6;  G = 1200; K = 0.50117210867155004; N = 0.8; O = 1.2; R = 0.50117210867155004 },
    @{ Row = 7;  G = 1200; K = 0.75175816300732523; N = 0.8; O = 1.2; R = 0.75175816300732523 },
    @{ Row = 8;  G = 1200; K = 0.10023442173431002; N = 0.8; O = 1.2; R = 0.10023442173431002 },
    @{ Row = 9;  G = 1200; K = 2.0604166666666663;  N = 0.8; O = 1.2; R = 2.0604166666666663 },
    @{ Row = 10; G = 3500; K = 1.4128571428571426;  N = 0.8; O = 1.2; R = 4.1208333333333327 },
    @{ Row = 11; G = 1200; K = 1.3300336730129598;  N = 0.8; O = 1.2; R = 1.3300336730129598 },
    @{ Row = 12; G = 3500; K = 1.4128571428571426;  N = 0.8; O = 1.2; R = 4.1208333333333327 },
    @{ Row = 13; G = 1200; K = 0.50117210867155004; N = 0.8; O = 1.2; R = 0.50117210867155004 },
    @{ Row = 14; G = 3500; K = 1.4128571428571426;  N = 0.8; O = 1.2; R = 4.1208333333333327 },
    @{ Row = 15; G = 3500; K = 2.4700625355954968;  N = 0.8; O = 1.2; R = 2.4700625355954968 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("G$n").Value = $r.G
    $ws.Range("G$n").NumberFormat = "General"

    $ws.Range("K$n").Value = $r.K
    $ws.Range("K$n").NumberFormat = "0.00"

    $ws.Range("N$n").Value = $r.N
    $ws.Range("O$n").Value = $r.O

    $ws.Range("R$n").Value = $r.R
    $ws.Range("R$n").NumberFormat = "0.00"
}

# New note row under the table
$ws.Range("K21").Value = " "

$ws.Range("I17").Select()
